# Generate Report for Handoff
# Rewrites the handoff-report identifiers (old GUID -> new GUID, new .xlf
# pair-hash, new handoff timestamps) on all three worksheets, while
# preserving each hyperlink's original target URL (only the visible
# display text / cell text changes), matching the authored diff.

$wb = $excel.ActiveWorkbook

$oldGuid = "508e0b21-cd50-48d2-9971-774d4a00fe5e"
$newGuid = "32294dc3-387f-4844-9e61-b6d684d4e3e4"

$oldHash = "1addd768f292c4ae4dcc1689f984b1ea628f2493"
$newHash = "0d243c4a6c9a539fbde1edef52be9f652ad55564"

$mdFileOld = "$oldGuid.md"
$mdFileNew = "$newGuid.md"

$zhFileOld = "$oldGuid.$oldHash.zh-cn.xlf"
$zhFileNew = "$newGuid.$newHash.zh-cn.xlf"

$deFileOld = "$oldGuid.$oldHash.de-de.xlf"
$deFileNew = "$newGuid.$newHash.de-de.xlf"

$zhTimeOld = "2016-03-04 08:46:44"
$zhTimeNew = "2016-03-04 08:47:33"

$deTimeOld = "2016-03-04 08:46:53"
$deTimeNew = "2016-03-04 08:47:43"

$configDisplay = ".localization-config"

$mdTargetUrl = "https://github.com/OpenLocalizationTest/oltest/blob/25c620c049b6f803578e6c7c686ad01676ca2318/e2e/$mdFileOld"
$configTargetUrl = "https://github.com/OpenLocalizationTest/oltest/blob/25c620c049b6f803578e6c7c686ad01676ca2318/.localization-config"
$zhXlfTargetUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/58563738b5f6a89225622f2e8ca5ea6610c32046/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$zhFileOld"
$deXlfTargetUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/74c9011b2c2d24b569857caf39dae4c1b1772c3a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$deFileOld"

# ---- Sheet 1: Overview ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A2").Value = $mdFileNew
$ws1.Range("A2").Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), $mdTargetUrl, "", "", $mdFileNew) | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), $configTargetUrl, "", "", $configDisplay) | Out-Null

# ---- Sheet 2: zh-cn ----
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A2").Value = $mdFileNew
$ws2.Range("C2").Value = $zhFileNew
$ws2.Range("D2").Value = $zhTimeNew
$ws2.Range("A2").Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $mdTargetUrl, "", "", $mdFileNew) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), $zhXlfTargetUrl, "", "", $zhFileNew) | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), $configTargetUrl, "", "", $configDisplay) | Out-Null

# ---- Sheet 3: de-de ----
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A2").Value = $mdFileNew
$ws3.Range("C2").Value = $deFileNew
$ws3.Range("D2").Value = $deTimeNew
$ws3.Range("A2").Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $mdTargetUrl, "", "", $mdFileNew) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), $deXlfTargetUrl, "", "", $deFileNew) | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), $configTargetUrl, "", "", $configDisplay) | Out-Null
